{"js": "const body = context.document.body;\n\n// 1. Title: \"STACK THE FLAGS ...\" -> \"STACK THE FLAG ...\" (drop the plural \"S\")\nconst flagsResults = body.search(\"FLAGS CAPTURE THE FLAG\", { matchCase: true });\nflagsResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < flagsResults.items.length; i++) {\n  flagsResults.items[i].insertText(\"FLAG CAPTURE THE FLAG\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 1b. Title: year \"2020\" -> \"2022\"\nconst yearResults = body.search(\"2020\", { matchCase: true });\nyearResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < yearResults.items.length; i++) {\n  yearResults.items[i].insertText(\"2022\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Disbursement paragraph: \"STACK the Flags' website\" -> \"STACK the Flag website\"\nconst apostropheResults = body.search(\"Flags' website\", { matchCase: true });\napostropheResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < apostropheResults.items.length; i++) {\n  apostropheResults.items[i].insertText(\"Flag website\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2b. Disbursement paragraph: drop the trailing sentence about physical\n// delivery of the C01Ns / residential address.\nconst tailSentence = \" For physical delivery of the GovTech STACK the Flags C01Ns, it will be sent to a single residential address.\";\nconst tailResults = body.search(tailSentence, { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tailResults.items.length; i++) {\n  tailResults.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Title: \"STACK THE FLAGS ...\" -> \"STACK THE FLAG ...\" (drop the plural \"S\")\n$find1 = $d.Content.Find\n$find1.Text = \"FLAGS CAPTURE THE FLAG\"\n$find1.Replacement.Text = \"FLAG CAPTURE THE FLAG\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)\n\n# 1b. Title: year \"2020\" -> \"2022\"\n$find2 = $d.Content.Find\n$find2.Text = \"2020\"\n$find2.Replacement.Text = \"2022\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n\n# 2. Disbursement paragraph: \"STACK the Flags' website\" -> \"STACK the Flag website\"\n$find3 = $d.Content.Find\n$find3.Text = \"Flags' website\"\n$find3.Replacement.Text = \"Flag website\"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2)\n\n# 2b. Disbursement paragraph: drop the trailing sentence about physical\n# delivery of the C01Ns / residential address.\n$find4 = $d.Content.Find\n$find4.Text = \" For physical delivery of the GovTech STACK the Flags C01Ns, it will be sent to a single residential address.\"\n$find4.Replacement.Text = \"\"\n$find4.Execute($find4.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find4.Replacement.Text, 2)\n"}
